# "danger zone gameplay content"
# Adds four new TCS food items (bbq pork, bbq chicken, broccoli salad, yogurt
# bar) right after the existing food list (inserted at row 38, pushing the
# old "TCS food description" / "thermometer" / "danger zone trivia" rows
# down by 4), and appends the danger-zone trivia question/answer strings for
# each new food item at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows right before the old row 38 (tcsFoodDesc_meats_title),
# shifting all subsequent rows down by 4.
$ws.Rows("38:41").Insert()

# New food key/value rows, each immediately followed (at the end of the
# sheet) by its danger-zone trivia block -- this mirrors the order the
# author actually typed the cells in (and therefore the resulting
# sharedStrings append order).
$ws.Range("A38").Value = "food_bbqPork"
$ws.Range("B38").Value = "Barbeque Pork"

$ws.Range("A76").Value = "dangerZoneTrivia2"
$ws.Range("B76").Value = "A tray of barbeque pork has been placed in the oven to be reheated. After 30 minutes, a server checks the internal temperature for 15 seconds. It reads 165° F."
$ws.Range("A77").Value = "dangerZoneTrivia20"
$ws.Range("B77").Value = "The temperature looks good, it is ready to be served."
$ws.Range("A78").Value = "dangerZoneTrivia21"
$ws.Range("B78").Value = "It still needs to be heated, put it back in the oven, and wait for another 30 minutes."
$ws.Range("A79").Value = "dangerZoneTrivia22"
$ws.Range("B79").Value = "It is not safe to be served, throw it out."

$ws.Range("A39").Value = "food_bbqChicken"
$ws.Range("B39").Value = "Barbeque Chicken"

$ws.Range("A80").Value = "dangerZoneTrivia3"
$ws.Range("B80").Value = "A server has been sent to check on the temperature of a tray of barbeque chicken. After checking the internal temperature for 15 seconds, it reads 127° F. The tray has been out for at least an hour."
$ws.Range("A81").Value = "dangerZoneTrivia30"
$ws.Range("B81").Value = "The tray needs to be heated up."
$ws.Range("A82").Value = "dangerZoneTrivia31"
$ws.Range("B82").Value = "The temperature looks good, leave it alone."

$ws.Range("A40").Value = "food_broccoliSalad"
$ws.Range("B40").Value = "Broccoli Salad"

$ws.Range("A83").Value = "dangerZoneTrivia4"
$ws.Range("B83").Value = "Several broccoli salad trays being served have been out for more than an hour. A server decided to check the temperature of each tray. They all read 40° F."
$ws.Range("A84").Value = "dangerZoneTrivia40"
$ws.Range("B84").Value = "The temperature looks good, leave it alone."
$ws.Range("A85").Value = "dangerZoneTrivia41"
$ws.Range("B85").Value = "These shouldn't be served, replace it with a new batch."
$ws.Range("A86").Value = "dangerZoneTrivia42"
$ws.Range("B86").Value = "It's too cold, they need to be heated up in a microwave."

$ws.Range("A41").Value = "food_yogurtBar"
$ws.Range("B41").Value = "Yogurt Bar"

$ws.Range("A87").Value = "dangerZoneTrivia5"
$ws.Range("B87").Value = "A server was requested to check on the yogurt batch being served at the yogurt bar. After checking the temperature, it reads 53° F. Looking at the time it was brought out, it's been several hours."
$ws.Range("A88").Value = "dangerZoneTrivia50"
$ws.Range("B88").Value = "This batch needs to be replaced."
$ws.Range("A89").Value = "dangerZoneTrivia51"
$ws.Range("B89").Value = "The temperature is fine, leave it as is."

# Match the final view state recorded in the saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 73
$ws.Range("A90").Select()

Write-Output "danger zone content added"
